$d = $word.ActiveDocument

# --- Step 1: merge the two runs around the old _GoBack bookmark in
# paragraph 3 into a single run, dropping the bookmark in the process.
# (Replace the *whole* paragraph body -- not just the first two runs --
# so the zero-width <w:proofErr/> marker that sits right at the old
# run/run boundary doesn't get reshuffled by the replace.)
$p3 = $d.Paragraphs.Item(3)
$mergeRange = $d.Range($p3.Range.Start, $p3.Range.End - 1)

$mergedFrag = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Я можу посилатися на додаток А чи додаток 1, але не можу писати додатки (а ДОДАТКИ можу!) або додаток </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>DiplomaAnalysis</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$mergeRange.InsertXML($mergedFrag)

# --- Step 2: append a brand-new paragraph after paragraph 3 with the
# "Наш дипломчик..." sentence (ru-RU runs) and move the _GoBack bookmark
# to the end of that new paragraph.
$p3 = $d.Paragraphs.Item(3)
$insertPoint = $p3.Range.End - 1
$insertRange = $d.Range($insertPoint, $insertPoint)

$newParaFrag = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve">Наш </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t>дипломчик</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>вирішує нашу проблему в тому, що я хочу отримати диплом, а так мені його не дадуть.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$insertRange.InsertXML($newParaFrag)
